$wb = $excel.ActiveWorkbook

# --- Echant1_impl: remove the lone "m+a" sample row3 breeding-status cell ---
$ws1 = $wb.Worksheets.Item("Echant1_impl")
$ws1.Activate()
$ws1.Range("C3").ClearContents()
[void]$ws1.Range("C3").Select()

# --- Params1_expl: drop the obsolete "Sylvia atricapilla / m+a / 10mn" block (rows 12-16) ---
$ws4 = $wb.Worksheets.Item("Params1_expl")
$ws4.Activate()
$ws4.Range("A12:G16").Delete()
[void]$ws4.Range("A12:XFD16").Select()

# --- Params2_expl: finish on this sheet (matches the saved activeTab/tabSelected state) ---
$ws5 = $wb.Worksheets.Item("Params2_expl")
$ws5.Activate()
[void]$ws5.Range("D22").Select()
